$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("wenden", "none", "none"),
    @("heben", "face/face022.jpg", "face"),
    @("drücken", "face/face024.jpg", "face"),
    @("rufen", "none", "none"),
    @("buchen", "house/house004.jpg", "house"),
    @("sehen", "house/house009.jpg", "house"),
    @("enden", "none", "none"),
    @("beißen", "face/face005.jpg", "face"),
    @("küssen", "house/house014.jpg", "house"),
    @("opfern", "none", "none"),
    @("lehnen", "house/house005.jpg", "house"),
    @("tragen", "face/face002.jpg", "face"),
    @("schalten", "none", "none"),
    @("pflegen", "house/house020.jpg", "house"),
    @("leuchten", "house/house001.jpg", "house"),
    @("klagen", "none", "none"),
    @("hoffen", "face/face028.jpg", "face"),
    @("münzen", "face/face026.jpg", "face"),
    @("sparen", "none", "none"),
    @("wandern", "house/house025.jpg", "house"),
    @("helfen", "house/house019.jpg", "house"),
    @("hören", "none", "none"),
    @("spielen", "face/face031.jpg", "face"),
    @("bieten", "face/face011.jpg", "face"),
    @("drohen", "none", "none"),
    @("geben", "face/face003.jpg", "face"),
    @("tanzen", "house/house008.jpg", "house"),
    @("orten", "none", "none"),
    @("schützen", "face/face018.jpg", "face"),
    @("treten", "house/house011.jpg", "house"),
    @("weigern", "none", "none"),
    @("planen", "house/house012.jpg", "house"),
    @("binden", "house/house018.jpg", "house"),
    @("ächzen", "none", "none"),
    @("kehren", "face/face001.jpg", "face"),
    @("dürfen", "face/face008.jpg", "face"),
    @("kosten", "none", "none"),
    @("schweben", "house/house000.jpg", "house"),
    @("wundern", "house/house028.jpg", "house"),
    @("dauern", "none", "none"),
    @("handeln", "face/face014.jpg", "face"),
    @("seufzen", "face/face009.jpg", "face"),
    @("stören", "none", "none"),
    @("knien", "face/face007.jpg", "face"),
    @("biegen", "face/face030.jpg", "face"),
    @("bremsen", "none", "none"),
    @("erben", "house/house007.jpg", "house"),
    @("schulden", "house/house002.jpg", "house")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
